# Auto refresh - 16-02-2026 13:33:53.19
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Excel_vs_ML")
$ws2 = $wb.Worksheets.Item("Feature_Importance")
$ws3 = $wb.Worksheets.Item("Exec_Summary")

# --- Sheet "Excel_vs_ML": update ML prediction columns (T,U,V,W) for refreshed rows ---
$ws1.Range("T5").Value = -18.51
$ws1.Range("U5").Value = 92.55000000000001
$ws1.Range("V5").Value = -94233.92318700001
$ws1.Range("W5").Value = "CRITICAL – Immediate Action"
$ws1.Range("T13").Value = -10.37
$ws1.Range("U13").Value = 51.84999999999999
$ws1.Range("V13").Value = -16758.44887
$ws1.Range("W13").Value = "MODERATE – Monitor Closely"
$ws1.Range("T15").Value = -9.32
$ws1.Range("U15").Value = 46.6
$ws1.Range("V15").Value = -39486.967612
$ws1.Range("W15").Value = "MODERATE – Monitor Closely"
$ws1.Range("T16").Value = -18.78
$ws1.Range("U16").Value = 93.9
$ws1.Range("V16").Value = -19897.610946
$ws1.Range("W16").Value = "CRITICAL – Immediate Action"
$ws1.Range("T18").Value = -10.8
$ws1.Range("U18").Value = 54
$ws1.Range("V18").Value = -30873.22884
$ws1.Range("W18").Value = "MODERATE – Monitor Closely"
$ws1.Range("T21").Value = 4.31
$ws1.Range("U21").Value = 21.55
$ws1.Range("V21").Value = 8481.104216
$ws1.Range("W21").Value = "LOW – Stable"
$ws1.Range("T23").Value = -10.51
$ws1.Range("U23").Value = 52.55
$ws1.Range("V23").Value = -37379.893326
$ws1.Range("W23").Value = "MODERATE – Monitor Closely"
$ws1.Range("T27").Value = -10.57
$ws1.Range("U27").Value = 52.84999999999999
$ws1.Range("V27").Value = -62569.166793
$ws1.Range("W27").Value = "MODERATE – Monitor Closely"
$ws1.Range("T28").Value = -17.56
$ws1.Range("U28").Value = 87.79999999999998
$ws1.Range("V28").Value = -68551.04583599999
$ws1.Range("W28").Value = "CRITICAL – Immediate Action"
$ws1.Range("T29").Value = -2.03
$ws1.Range("U29").Value = 10.15
$ws1.Range("V29").Value = -2603.856437
$ws1.Range("W29").Value = "LOW – Stable"
$ws1.Range("T30").Value = -10.39
$ws1.Range("U30").Value = 51.95000000000001
$ws1.Range("V30").Value = -55485.14555
$ws1.Range("W30").Value = "MODERATE – Monitor Closely"
$ws1.Range("T32").Value = 3.1
$ws1.Range("U32").Value = 15.5
$ws1.Range("V32").Value = 18487.74683
$ws1.Range("W32").Value = "LOW – Stable"
$ws1.Range("T34").Value = -17.92
$ws1.Range("U34").Value = 89.60000000000001
$ws1.Range("V34").Value = -102908.723456
$ws1.Range("W34").Value = "CRITICAL – Immediate Action"
$ws1.Range("T36").Value = -21.28
$ws1.Range("U36").Value = 100
$ws1.Range("V36").Value = -75461.222928
$ws1.Range("W36").Value = "CRITICAL – Immediate Action"

# Remove the Early_Warning column (X) - no longer produced by the refreshed pipeline
$ws1.Columns.Item(24).Delete()

# --- Sheet "Feature_Importance": replace with refreshed feature ranking ---
$ws2.Range("A2").Value = "Budget_Pct"
$ws2.Range("B2").Value = 0.4243173945580014
$ws2.Range("A3").Value = "Gap_Pct"
$ws2.Range("B3").Value = 0.4213112563496177
$ws2.Range("A4").Value = "Acceleration"
$ws2.Range("B4").Value = 0.06059531511456343
$ws2.Range("A5").Value = "Spend_Velocity"
$ws2.Range("B5").Value = 0.02982975275704113
$ws2.Range("A6").Value = "DSP_enc"
$ws2.Range("B6").Value = 0.0185569584511927
$ws2.Range("A7").Value = "Total_Budget"
$ws2.Range("B7").Value = 0.0134733367276791
$ws2.Range("A8").Value = "Time_Pct"
$ws2.Range("B8").Value = 0.01223233831247451
$ws2.Range("A9").Value = "Spend_to_Date"
$ws2.Range("B9").Value = 0.01064231977620999
$ws2.Range("A10").Value = "Days_Remaining"
$ws2.Range("B10").Value = 0.003280289219582662
$ws2.Range("A11").Value = "Flight_Days"
$ws2.Range("B11").Value = 0.00323182397208861
$ws2.Range("A12").Value = "Days_Elapsed"
$ws2.Range("B12").Value = 0.002529214761548717

# --- Sheet "Exec_Summary": refresh summary metrics + timestamp ---
$ws3.Range("B2").Value = 0.0463
$ws3.Range("B3").Value = -579240.38
$ws3.Range("B4").Value = "2026-02-16 08:03 UTC"

